$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2 through 13, column C: update the "Förändrad" (changed/expiry) date
# from 2023-09-16 (serial 45185) to 2023-10-05 (serial 45204).
for ($row = 2; $row -le 13; $row++) {
    $ws.Cells.Item($row, 3).Value = 45204
}
